$wb = $excel.ActiveWorkbook

# --- Sheet "Artistes" ---
$wsArtistes = $wb.Worksheets.Item("Artistes")

# Put each artist on the right scene number (column F)
$wsArtistes.Range("F3").Value = 2
$wsArtistes.Range("F4").Value = 3
$wsArtistes.Range("F5").Value = 4
$wsArtistes.Range("F6").Value = 5
$wsArtistes.Range("F7").Value = 6

# --- Sheet "Scenes" ---
$wsScenes = $wb.Worksheets.Item("Scenes")

# Fix scene/artist-count column (B) -> every scene now has 1 artist
$wsScenes.Range("B3").Value = 1
$wsScenes.Range("B4").Value = 1
$wsScenes.Range("B5").Value = 1
$wsScenes.Range("B6").Value = 1
$wsScenes.Range("B7").Value = 1

# Re-apply formatting on column F (Numéro de scène) to match column E's style
$wsScenes.Range("F2:F7").Style = "Normal"

# --- Restore cursor positions on each sheet (Artistes selected first, then
#     Scenes last, so Scenes remains the active/visible tab, like before) ---
$wsArtistes.Range("A2").Select()

$wsScenes.Activate()
$wsScenes.Range("B10").Select()
